$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 6 additions (Is computed on the fly?) ---
$ws.Range("C6").Value = "no"
$ws.Range("E6").Value = "no"
$ws.Range("H6").Value = "no"
$ws.Range("I6").Value = "no"
$ws.Range("J6").Value = "no"
$ws.Range("K6").Value = "no"
$ws.Range("L6").Value = "yes"
$ws.Range("M6").Value = "yes"
$ws.Range("N6").Value = "no"
$ws.Range("O6").Value = "no def 1.0"
$ws.Range("P6").Value = "yes"

# --- Row 7 (label change + new data) ---
$ws.Range("A7").Value = "May be in constructor"
$ws.Range("C7").Value = "yes"
$ws.Range("D7").Value = "yes"
$ws.Range("E7").Value = "y"
$ws.Range("G7").Value = "yes"
$ws.Range("H7").Value = "yes"
$ws.Range("I7").Value = "no"
$ws.Range("J7").Value = "yes"
$ws.Range("K7").Value = "no"
$ws.Range("L7").Value = "yes"
$ws.Range("M7").Value = "yes"
$ws.Range("N7").Value = "yes"
$ws.Range("O7").Value = "yes"
$ws.Range("P7").Value = "yes"

# --- Row 8 (Source class row) ---
$ws.Range("D8").Value = "Path"
$ws.Range("E8").Value = "Graphic"

# --- New rows 12 and 13 (new constructor test rows) ---
$ws.Range("A12").Value = "Arc_CreateArc_ctor3DeflectingLeft90_IsCorrect"
$ws.Range("B12").Value = "n"
$ws.Range("A13").Value = "Arc_CreateArc_ctor3DeflectingRight90_IsCorrect"
$ws.Range("B13").Value = "n"

# --- Column widths (target stored widths: col A = 61.6640625, col B = 11.5546875) ---
$ws.Columns.Item(1).ColumnWidth = 60.833333333333336
$ws.Columns.Item(2).ColumnWidth = 10.666666666666666

# --- Selection change (bottomRight pane active cell H6 -> C6) ---
$ws.Range("C6").Select()
